# Actualización automática 2025-07-31 13:10:08
# Updates "julio" (July) sales figures across the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-client / per-category detail
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# MADECOR-HOME CENTER S.A.S. (row 29): FREGADEROS DE COCINA + LAVABOS sales recorded
$wsGrupo.Range("E29").Value = 166.95
$wsGrupo.Range("I29").Value = 313.2

# MARQUEZ CELI ALFREDO NICANOR (row 30): NO RESURTIBLES sale recorded
$wsGrupo.Range("P30").Value = 25.85

# RUIZ TINIZARAY YOHANNA MARYURI (row 49): 240X80 PORCELANATO sale recorded
$wsGrupo.Range("D49").Value = 950.4

# Row 56 totals — "N de 54" completion counters, +1 for each column that
# just went from 0 to a real value above.
$wsGrupo.Range("D56").Value = "2 de 54"
$wsGrupo.Range("E56").Value = "5 de 54"
$wsGrupo.Range("I56").Value = "7 de 54"
$wsGrupo.Range("P56").Value = "2 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" — per-client monthly totals ("julio" = column F)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F29").Value = 4226.2
$wsMensual.Range("F30").Value = 25.85
$wsMensual.Range("F49").Value = 1971.65
$wsMensual.Range("F56").Value = 84614.60000000001

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — per-category compliance summary
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 1408.32
$wsCumpl.Range("E3").Value = 26048.6876
$wsCumpl.Range("F3").Value = 0.05129182394952609

# Row 4: FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 1179.68
$wsCumpl.Range("E4").Value = -176.6800000000001
$wsCumpl.Range("F4").Value = 1.176151545363908

# Row 8: LAVABOS
$wsCumpl.Range("D8").Value = 1006.06
$wsCumpl.Range("E8").Value = -6.059999999999945
$wsCumpl.Range("F8").Value = 1.00606

# Row 10: NO RESURTIBLES
$wsCumpl.Range("D10").Value = 143.94
$wsCumpl.Range("E10").Value = 1156.56
$wsCumpl.Range("F10").Value = 0.1106805074971165

# Row 19: TOTAL
$wsCumpl.Range("D19").Value = 84614.60000000001
$wsCumpl.Range("E19").Value = 29091.85064517915
$wsCumpl.Range("F19").Value = 0.7441495141207052
